$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 69, pushing existing rows 69:167 down to 70:168.
$ws.Rows(69).Insert()

# Populate the newly inserted row 69 with the new daily price entry.
$ws.Range("A69").Value = 4
$ws.Range("B69").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C69").Value = "Los Lagos"
$ws.Range("D69").Value = 44467
$ws.Range("E69").Value = 10
$ws.Range("F69").Value = 100112040
$ws.Range("G69").Value = "Cilantro"
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 300
$ws.Range("K69").Value = 12000
$ws.Range("L69").Value = 12000
$ws.Range("M69").Value = 12000
$ws.Range("N69").Value = "$/caja 36 atados"
$ws.Range("O69").Value = "Región Metropolitana"
$ws.Range("P69").Value = 333
$ws.Range("Q69").Value = 36
$ws.Range("R69").Value = "Hortaliza"
